$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted as row 32; every existing row from
# 32 onward (old rows 32-65) shifts down by one (to 33-66).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 44571
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100101
$ws.Range("H32").Value = "Berries"
$ws.Range("I32").Value = 100101008
$ws.Range("J32").Value = "Mora"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 6000
$ws.Range("O32").Value = 6000
$ws.Range("P32").Value = 6000
$ws.Range("Q32").Value = "`$/bandeja 2 kilos"
$ws.Range("R32").Value = "Provincia de Linares"
$ws.Range("S32").Value = 3000
$ws.Range("T32").Value = 2
